$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 16: self-assessed score changes from 70 to 91
$ws.Range("M16").Value = 91

# Row 18: leader-approval rating changes, score, and approving leader name
$ws.Range("L18").Value = "Không đảm bảo chất lượng"
$ws.Range("M18").Value = 40
$ws.Range("N18").Value = "Chan Dan"

# Row 19: self-assessed score changes from 82 to 95
$ws.Range("M19").Value = 95
